# Generate Report for Handoff
# Adds a new tracked file (8afa32e3-faa1-4036-8e5b-c5eda857ff25.md) as row 9
# on the Overview / zh-cn / de-de sheets, mirroring the existing rows.

$wb = $excel.ActiveWorkbook

$fileId   = "8afa32e3-faa1-4036-8e5b-c5eda857ff25"
$mdName   = "$fileId.md"
$zhHash   = "d03c38702865a5f22fd38d05a92e06aa5c697f06"
$zhXlf    = "$fileId.$zhHash.zh-cn.xlf"
$deXlf    = "$fileId.$zhHash.de-de.xlf"

$mdUrl      = "https://github.com/OpenLocalizationTest/oltest/blob/c13e91bfdbc08260c3d9c60df3c0886866b433a7/e2e/$mdName"
$zhXlfUrl   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ca1cda8444eb29c011fb3836247aa8ad43a5c0ee/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf"
$deXlfUrl   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ca638b4bf88f531f944219eec792bb7029b4eb82/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf"

$readyStatus = "Ready for handoff"

# ---------------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")
$row = 9

$ws1.Cells.Item($row, 1).Value = $mdName
$ws1.Cells.Item($row, 2).Value = $readyStatus
$ws1.Cells.Item($row, 3).Value = $readyStatus
$ws1.Cells.Item($row, 4).Value = "2016-03-23 10:46:27"

$ws1.Hyperlinks.Add($ws1.Cells.Item($row, 1), $mdUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdName) | Out-Null

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Cells.Item($row, 1).Value = $mdName
$ws2.Cells.Item($row, 2).Value = ".md"
$ws2.Cells.Item($row, 3).Value = $readyStatus
$ws2.Cells.Item($row, 4).Value = $zhXlf
$ws2.Cells.Item($row, 5).Value = "2016-03-23 10:46:24"
$ws2.Cells.Item($row, 8).Value = "0001-01-01 00:00:00"
$ws2.Cells.Item($row, 10).Value = "Include"

$ws2.Hyperlinks.Add($ws2.Cells.Item($row, 1), $mdUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdName) | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item($row, 4), $zhXlfUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $zhXlf) | Out-Null

# ---------------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Cells.Item($row, 1).Value = $mdName
$ws3.Cells.Item($row, 2).Value = ".md"
$ws3.Cells.Item($row, 3).Value = $readyStatus
$ws3.Cells.Item($row, 4).Value = $deXlf
$ws3.Cells.Item($row, 5).Value = "2016-03-23 10:46:27"
$ws3.Cells.Item($row, 8).Value = "0001-01-01 00:00:00"
$ws3.Cells.Item($row, 10).Value = "Include"

$ws3.Hyperlinks.Add($ws3.Cells.Item($row, 1), $mdUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdName) | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item($row, 4), $deXlfUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $deXlf) | Out-Null

Write-Host "Row 9 added to Overview, zh-cn, de-de sheets."
